$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab11")

# Row 67: only column I changes
$ws.Range("I67").Value = 62.9267914160414

# Row 97: "Afrique, Etats fragiles" - all data columns C:J change
$ws.Range("C97").Value = 0.563188849578700
$ws.Range("D97").Value = 26.5373753562805
$ws.Range("E97").Value = 76.4741122245553
$ws.Range("F97").Value = 0.735761102116110
$ws.Range("G97").Value = 28.5881777963223
$ws.Range("H97").Value = 842905.229565056
$ws.Range("I97").Value = 57.4983111689084
$ws.Range("J97").Value = 25.267971338377

# Row 98: "RDM, Etats fragiles" - all data columns C:J change
$ws.Range("C98").Value = 4.15906858058393
$ws.Range("D98").Value = 32.1352928495309
$ws.Range("E98").Value = 95.5257154881166
$ws.Range("F98").Value = 6.87584224912458
$ws.Range("G98").Value = 50.9000828089096
$ws.Range("H98").Value = 1792876.98733909
$ws.Range("I98").Value = 80.4138105601341
$ws.Range("J98").Value = 29.6264260164983

$wb.Save()
